$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("A18").Value = "billion 2020 dollars"
$ws.Range("A21").Value = "million 2020 dollars"
$ws.Range("A24").Value = "2020 dollars"
$ws.Range("B29").Value = 'which in this case is "2012 dollars per 2020 dollar."'
$ws.Range("B30").Value = "2012 dollars are worth more than 2020 dollars, so we need a"

$ws.Range("A26").Value = 0.88711067149387013
